$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "F1a"
$ws.Range("A5").Value = "F1b"
$ws.Range("A6").Value = "F2a"
$ws.Range("A7").Value = "F2b"

$ws.Range("F17").Select()
